# Apply updated dSF (column F) values pulled from a repull of source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = -5
$ws.Range("F3").Value = -6
$ws.Range("F6").Value = -5
$ws.Range("F11").Value = 0
